# Capstone 3, Presentation.pptx - apply edits described in the commit
# (upload dated 2023-04-19) on top of the 2023-04-17 version.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 1 ("Title Slide"): subtitle date line "March 2022" -> "March 2023"
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.Name -eq "Subtitle 2") {
        $tr = $shp.TextFrame.TextRange
        $paras = $tr.Paragraphs()
        for ($pi = 1; $pi -le $paras.Count; $pi++) {
            $para = $tr.Paragraphs($pi, 1)
            if ($para.Text -eq "March 2022") {
                $para.Text = "March 2023"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Footer "Fixed" date placeholders: "4/3/2023" -> "4/18/2023"
#    (lives on the slide master and every custom layout)
# ---------------------------------------------------------------------------
$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "*Date*" -and $shp.TextFrame.TextRange.Text -eq "4/3/2023") {
        $shp.TextFrame.TextRange.Text = "4/18/2023"
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "*Date*" -and $shp.TextFrame.TextRange.Text -eq "4/3/2023") {
            $shp.TextFrame.TextRange.Text = "4/18/2023"
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Slide 27 ("Next steps / future work"): append a new bullet paragraph
#    to the "Content Placeholder 2" shape, with a mix of regular and
#    italic runs.
# ---------------------------------------------------------------------------
$slide27 = $p.Slides.Item(27)
for ($i = 1; $i -le $slide27.Shapes.Count; $i++) {
    $shp = $slide27.Shapes.Item($i)
    if ($shp.Name -eq "Content Placeholder 2") {

        $tr = $shp.TextFrame.TextRange
        $lastIdx = $tr.Paragraphs().Count
        $lastPara = $tr.Paragraphs($lastIdx, 1)

        $run1 = "We had a high percentage of success among lines "
        $run2 = "already determined to be likely to succeed"
        $run3 = " " + [char]0x2013 + " consider analyzing more lines in the original data for a more accurate "
        $run4 = "overall picture"
        $fullText = $run1 + $run2 + $run3 + $run4

        $lastPara.InsertAfter("`r" + $fullText)

        $tr2 = $shp.TextFrame.TextRange
        $newIdx = $lastIdx + 1
        $newPara = $tr2.Paragraphs($newIdx, 1)
        # Match the outline level of the earlier top-level bullets (no
        # sub-level indent) rather than inheriting level 2 from the
        # preceding "Focus analysis..." bullet.
        $newPara.IndentLevel = 1

        $start = 1
        $sub2 = $newPara.Characters($start + $run1.Length, $run2.Length)
        $sub2.Font.Italic = $true
    }
}

Write-Output "edits applied"
